# Applies the "Proyecto final POO v2.0" edits:
#  1) Paragraph "Utilizar templates ..." - merge the "Utilizar " / "templates"
#     (proof-err wrapped) / " para permitir ..." runs into a single run, and
#     split off the trailing " Opcional" into its own space-run + "Opcional"
#     run (no leading space).
#  2) Paragraph "Uso de MVC en el proyecto" - append " Opcional" (highlighted
#     green, "Opcional" also colored red) as new runs.

$d = $word.ActiveDocument

# Locate the two target paragraphs by their current text.
$templatesPara = $null
$mvcPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Utilizar templates*") {
        $templatesPara = $p
    }
    if ($t -like "Uso de MVC en el proyecto*") {
        $mvcPara = $p
    }
}

if ($templatesPara -eq $null) {
    throw "Could not find the 'Utilizar templates ...' paragraph"
}
if ($mvcPara -eq $null) {
    throw "Could not find the 'Uso de MVC en el proyecto' paragraph"
}

# --- 1) "Utilizar templates ..." paragraph -------------------------------
$r1 = $templatesPara.Range
$body1 = $d.Range($r1.Start, $r1.End - 1)

$templatesInner = '<w:r w:rsidRPr="00461CE2"><w:rPr><w:highlight w:val="red"/></w:rPr>' +
    '<w:t>Utilizar templates para permitir el almacenamiento de diferentes tipos de animales y visitantes en listas din' + [char]0x00E1 + 'micas.</w:t></w:r>' +
    '<w:r w:rsidRPr="00461CE2"><w:rPr><w:highlight w:val="red"/></w:rPr>' +
    '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidR="00461CE2"><w:rPr><w:highlight w:val="red"/></w:rPr>' +
    '<w:t>Opcional</w:t></w:r>'

$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $templatesInner + '</w:p>'
$body1.InsertXML($xml1)

# --- 2) "Uso de MVC en el proyecto" paragraph -----------------------------
$r2 = $mvcPara.Range
$body2 = $d.Range($r2.Start, $r2.End - 1)

$mvcInner = '<w:r w:rsidRPr="00461CE2"><w:rPr><w:highlight w:val="green"/></w:rPr>' +
    '<w:t>Uso de MVC en el proyecto</w:t></w:r>' +
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr>' +
    '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:highlight w:val="green"/></w:rPr>' +
    '<w:t>O</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:highlight w:val="green"/></w:rPr>' +
    '<w:t>pcional</w:t></w:r>'

$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $mvcInner + '</w:p>'
$body2.InsertXML($xml2)

Write-Output "Done"
